$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date-serial value that was bumped by one day
# (46061 -> 46062) for every data row (rows 2 through 388).
$ws.Range("C2:C388").Value = 46062
